$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value happens to look like a plain number to
# Excel's auto-detection (e.g. "308.79") must stay text, matching the
# original inlineStr cells in the source file. Temporarily mark each such
# cell as Text before writing the value, then restore the default (no
# explicit format) style afterwards so no stray formatting is introduced.

$ws.Range('D2').Value = '23.880.56'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.645.76'
$ws.Range('E3').Value = '  +1.82%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '308.79'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('E7').Value = '  -0.69%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3821'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +3.93%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.347'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('E11').Value = '  +0.12%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08435'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '23.82'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.084'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.759'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.08%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001308'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.68%  '
$ws.Range('D17').Value = '1.646.18'
$ws.Range('E17').Value = '  +1.64%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '94.36'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06974'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.21%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '19.66'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.854'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.64%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '13.54'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('D24').Value = '23.869.64'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.479'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.008'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +5.55%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '22.01'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.50%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '152.66'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.76%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.425'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.46%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '138.70'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.731'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.489'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('D33').Value = '1.827.88'
$ws.Range('E33').Value = '  +1.97%  '
$ws.Range('E34').Value = '  +4.85%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.08028'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02945'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.66%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.677'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '10.85'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.42%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2669'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.09091'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.7511'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '13.38'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.411'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.47%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.15'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6889'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.435'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.69%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.065'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.08264'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '133.95'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('E51').Value = '  +0.93%  '
